# Tuntiseuranta.xlsx update
# - Jaana sheet: add "Tunnit:" summary row (row 5), add a new work-log entry (44992 / 4h / Servicelistin...)
# - Jarno sheet: add "Tunnit:" summary row (row 5)
# - Maarit sheet: unchanged

$wb = $excel.ActiveWorkbook

$wsJaana = $wb.Worksheets.Item("Jaana")
$wsJarno = $wb.Worksheets.Item("Jarno")
$wsMaarit = $wb.Worksheets.Item("Maarit")

# ---------------------------------------------------------------------------
# Jaana sheet
# ---------------------------------------------------------------------------

# Insert a new row above the current row 5 (blank separator row) for the
# "Tunnit:" running-total line, pushing everything else down by one row.
$wsJaana.Rows.Item(5).Insert()

# Copy formatting of the existing header-label cell (A4) onto the new label
# cell so it keeps the same look (bold-ish label style used for "Ryhma:" /
# "Henkilo:" rows).
$wsJaana.Range("A4").Copy()
$wsJaana.Range("A5").PasteSpecial(-4122)

$wsJaana.Range("A5").Value2 = "Tunnit: "
$wsJaana.Range("B5").Formula = "=B39"

# New work entry that used to be an empty row in the table (old row 17,
# shifted to row 18 after the insert above).
$wsJaana.Range("A8").Copy()
$wsJaana.Range("A18").PasteSpecial(-4122)
$wsJaana.Range("A18").Value2 = 44992
$wsJaana.Range("B18").Value2 = 4
$wsJaana.Range("C18").Value2 = "Servicelistin, Customerlistin ja Invoicelistin tietojen kuljetus json-tiedostojen kautta ikkunoiden välillä."
$wsJaana.Rows.Item(18).RowHeight = 48

$wsJaana.Range("A6").Select()

# ---------------------------------------------------------------------------
# Jarno sheet
# ---------------------------------------------------------------------------

$wsJarno.Rows.Item(5).Insert()

$wsJarno.Range("A5").Value2 = "Tunnit:"
$wsJarno.Range("B5").Formula = "=B39"

$wsJarno.Range("A4").Copy()
$wsJarno.Range("A5").PasteSpecial(-4122)

$wsJarno.Range("A6").Select()

$wb.Save()
